$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '56.516.53'
$ws.Cells.Item(2,5).Value = '  -2.41%  '

$ws.Cells.Item(3,4).Value = '2.385.58'
$ws.Cells.Item(3,5).Value = '  -2.82%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '0.999'
$ws.Cells.Item(4,5).Value = '  -0.07%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '505.23'
$ws.Cells.Item(5,5).Value = '  -4.20%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '131.45'
$ws.Cells.Item(6,5).Value = '  +0.48%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.996'
$ws.Cells.Item(7,5).Value = '  -0.49%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.555'
$ws.Cells.Item(8,5).Value = '  -1.51%  '

$ws.Cells.Item(9,4).Value = '2.405.32'
$ws.Cells.Item(9,5).Value = '  -2.01%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.0964'
$ws.Cells.Item(10,5).Value = '  -1.70%  '

$ws.Cells.Item(11,5).Value = '  -0.91%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.319'
$ws.Cells.Item(12,5).Value = '  -0.85%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '4.65'
$ws.Cells.Item(13,5).Value = '  -5.89%  '

$ws.Cells.Item(14,4).Value = '2.810.77'
$ws.Cells.Item(14,5).Value = '  -2.77%  '

$ws.Cells.Item(15,4).Value = '56.381.22'
$ws.Cells.Item(15,5).Value = '  -2.49%  '

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '21.62'
$ws.Cells.Item(16,5).Value = '  -0.89%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.0000132'
$ws.Cells.Item(17,5).Value = '  -0.79%  '

$ws.Cells.Item(18,4).Value = '2.399.90'
$ws.Cells.Item(18,5).Value = '  -2.12%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '10.17'
$ws.Cells.Item(19,5).Value = '  -1.87%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '310.99'
$ws.Cells.Item(20,5).Value = '  -1.80%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '4.04'
$ws.Cells.Item(21,5).Value = '  -2.83%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '6.34'
$ws.Cells.Item(22,5).Value = '  +4.46%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '0.998'
$ws.Cells.Item(23,5).Value = '  -0.17%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '5.59'
$ws.Cells.Item(24,5).Value = '  -4.44%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '64.46'
$ws.Cells.Item(25,5).Value = '  -0.81%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '0.994'
$ws.Cells.Item(26,5).Value = '  -0.76%  '

$ws.Cells.Item(27,4).Value = '2.487.63'
$ws.Cells.Item(27,5).Value = '  -3.28%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '0.378'
$ws.Cells.Item(28,5).Value = '  -7.36%  '

$ws.Cells.Item(29,5).Value = '  -4.89%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '7.38'
$ws.Cells.Item(30,5).Value = '  +1.25%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '172.61'
$ws.Cells.Item(31,5).Value = '  +0.27%  '

$ws.Cells.Item(32,4).Value = '0.0₃0722'
$ws.Cells.Item(32,5).Value = '  -1.83%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '1.67'
$ws.Cells.Item(33,5).Value = '  -1.48%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '6.13'
$ws.Cells.Item(34,5).Value = '  -0.03%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.13'
$ws.Cells.Item(35,5).Value = '  -2.81%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.993'
$ws.Cells.Item(37,5).Value = '  -0.65%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '17.92'
$ws.Cells.Item(38,5).Value = '  +0.32%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '1.22'
$ws.Cells.Item(39,5).Value = '  +3.11%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '3.80'
$ws.Cells.Item(40,5).Value = '  -0.27%  '

$ws.Cells.Item(41,2).Value = 'SuiNetwork'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.807'
$ws.Cells.Item(41,5).Value = '  +0.26%  '

$ws.Cells.Item(42,2).Value = 'OKB'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '35.94'
$ws.Cells.Item(42,5).Value = '  -0.71%  '

$ws.Cells.Item(43,5).Value = '  -1.48%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '129.60'
$ws.Cells.Item(44,5).Value = '  +3.61%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '3.37'
$ws.Cells.Item(45,5).Value = '  -1.36%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '4.92'
$ws.Cells.Item(46,5).Value = '  +1.42%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '254.18'
$ws.Cells.Item(47,5).Value = '  -5.13%  '

$ws.Cells.Item(48,2).Value = 'Mantle'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.562'
$ws.Cells.Item(48,5).Value = '  -3.84%  '

$ws.Cells.Item(49,2).Value = 'Stellar'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '0.0907'
$ws.Cells.Item(49,5).Value = '  -2.49%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.0490'
$ws.Cells.Item(50,5).Value = '  -1.18%  '

$ws.Cells.Item(51,2).Value = 'VeChain'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.0210'
$ws.Cells.Item(51,5).Value = '  -0.64%  '
